$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.349.15"
$ws.Range("E2").Value = "  -5.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.458.27"
$ws.Range("E3").Value = "  -8.36%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.40"
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.86"
$ws.Range("E6").Value = "  -6.78%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.477.49"
$ws.Range("E9").Value = "  -7.76%  "
$ws.Range("E10").Value = "  -5.99%  "
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.32"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.898.88"
$ws.Range("E14").Value = "  -8.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.16"
$ws.Range("E15").Value = "  -8.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.273.00"
$ws.Range("E16").Value = "  -5.70%  "
$ws.Range("E17").Value = "  -5.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.517.94"
$ws.Range("E18").Value = "  -6.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("E19").Value = "  -6.13%  "
$ws.Range("E20").Value = "  -5.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.52"
$ws.Range("E21").Value = "  -5.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.968"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("E23").Value = "  -9.05%  "
$ws.Range("E24").Value = "  -9.16%  "
$ws.Range("E25").Value = "  -4.31%  "
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.28"
$ws.Range("E29").Value = "  -7.00%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.78"
$ws.Range("E30").Value = "  -6.66%  "
$ws.Range("E31").Value = "  -6.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0773"
$ws.Range("E32").Value = "  -10.04%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.59"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.58"
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("E36").Value = "  -7.39%  "
$ws.Range("E37").Value = "  -5.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.97"
$ws.Range("E39").Value = "  -5.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "321.07"
$ws.Range("E40").Value = "  -8.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.77"
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("E42").Value = "  -12.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.70"
$ws.Range("E43").Value = "  -7.49%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.73"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0525"
$ws.Range("E48").Value = "  -6.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.05"
$ws.Range("E49").Value = "  -8.50%  "
$ws.Range("E50").Value = "  -8.68%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  -5.16%  "
